$d = $word.ActiveDocument

# The document ends with a paragraph containing a single space.
# We need to insert a new "Accountant" section (bold heading + 4 items +
# a blank paragraph) right before that trailing paragraph.
$lastParaIndex = $d.Paragraphs.Count
$insertRange = $d.Paragraphs.Item($lastParaIndex).Range
$insertRange.Collapse(1)   # wdCollapseStart

$newText = "Accountant`rItem Purchase Manage`rDoctor Salary`rEmployee Salary`rSalary history`r`r"
$insertRange.InsertBefore($newText)

# Make the "Accountant" heading paragraph bold, matching the style used
# by the other section headings (Admin, Doctor, Patient, Pharmacist).
$headingPara = $d.Paragraphs.Item($lastParaIndex)
$headingPara.Range.Font.Bold = 1
